$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Exercise" / " 2" -> single run "Exercise 2", and drop the old
#    _GoBack bookmark that currently sits right after it.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Exercise 2", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Exercise 2", 2)

$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) "Lab 2 ..." -> "Lab 3 ..." while preserving the original run
#    layout (six runs) of that paragraph, then re-insert the
#    _GoBack bookmark right after the run that now reads "3".
# ------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("Lab 2")
$startLab = $idx
$pos2 = $idx + 4        # offset of the "2" character
$pos2End = $pos2 + 1    # offset right after the "2" character

# Remember the paragraph's original run-boundary offsets (relative to
# the whole story) so we can restore them after the text edit, which
# otherwise collapses every run in the paragraph into one.
$boundaries = @($pos2, $pos2End, $pos2End + 1, $pos2End + 2, $pos2End + 3)

# Change the digit itself.
$rngDigit = $d.Range($pos2, $pos2End)
$rngDigit.Text = "3"

# Force the paragraph's runs back apart at their original boundaries
# by dropping (and immediately removing) temporary bookmarks there -
# inserting a bookmark splits the run it lands in without merging
# anything back together.
$i = 0
foreach ($off in $boundaries) {
    $tmpName = "zzTempSplit" + $i
    $tmpRng = $d.Range($off, $off)
    $d.Bookmarks.Add($tmpName, $tmpRng)
    $i = $i + 1
}
$i = 0
foreach ($off in $boundaries) {
    $tmpName = "zzTempSplit" + $i
    $d.Bookmarks($tmpName).Delete()
    $i = $i + 1
}

# Finally, place the _GoBack bookmark right after the "3".
$rngMark = $d.Range($pos2End, $pos2End)
$d.Bookmarks.Add("_GoBack", $rngMark)
